# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K" in row 1) held a running/derived
# "Strike#" count. Re-generating the save data recomputes this column from
# the authoritative per-game knight ("K") counts, so each row's K value is
# rewritten with its freshly-calculated figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row -> recalculated K value
$kValues = [ordered]@{
    2  = 0
    3  = 0
    4  = 0
    5  = 4
    6  = 1
    7  = 4
    8  = 7
    9  = 2
    10 = 5
    11 = 7
    12 = 4
    13 = 4
    14 = 6
    15 = 4
    16 = 16
    17 = 6
    18 = 11
    19 = 4
    20 = 1
    21 = 4
    22 = 4
    23 = 4
    24 = 10
    25 = 2
    26 = 3
    27 = 3
    28 = 1
    29 = 6
    30 = 3
    31 = 6
    32 = 5
    33 = 3
    34 = 1
    35 = 1
    36 = 5
    37 = 1
    38 = 2
    39 = 0
    40 = 1
    41 = 0
    42 = 0
    43 = 4
    44 = 6
    45 = 2
    46 = 5
    47 = 2
}

$col = 7  # column G = "K"

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, $col).Value = $kValues[$row]
}
